$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 6376475.5
$ws.Range("I98").Value = 9316595
$ws.Range("J98").Value = 6216.6665
$ws.Range("K98").Value = 9316595
$ws.Range("L98").Value = 6216.6665
$ws.Range("M98").Value = -9315097
$ws.Range("N98").Value = -9212.666499999999
$ws.Range("H122").Value = 6376475.5
$ws.Range("I122").Value = 9316595
$ws.Range("J122").Value = 6216.6665
$ws.Range("K122").Value = 27949785
$ws.Range("L122").Value = 18649.9995
$ws.Range("M122").Value = -27947335
$ws.Range("N122").Value = -23549.9995
$ws.Range("H129").Value = 942.61536
$ws.Range("I129").Value = 324.54544
$ws.Range("J129").Value = 1068.5186
$ws.Range("K129").Value = 973.63632
$ws.Range("L129").Value = 3205.5558
$ws.Range("M129").Value = 4026.36368
$ws.Range("N129").Value = -13205.5558
$ws.Range("H132").Value = 4040.9412
$ws.Range("I132").Value = 4255.8965
$ws.Range("J132").Value = 2794.2
$ws.Range("K132").Value = 12767.6895
$ws.Range("L132").Value = 8382.599999999999
$ws.Range("M132").Value = -10237.6895
$ws.Range("N132").Value = -13442.6
$ws.Range("H138").Value = 3513.24
$ws.Range("I138").Value = 1438.3334
$ws.Range("J138").Value = 5428.5386
$ws.Range("K138").Value = 4315.0002
$ws.Range("L138").Value = 16285.6158
$ws.Range("M138").Value = 824.9997999999996
$ws.Range("N138").Value = -26565.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2565760.5
$ws.Range("I61").Value = 3510145.8
$ws.Range("J61").Value = 2428.5715
$ws.Range("K61").Value = 3510145.8
$ws.Range("L61").Value = 2428.5715
$ws.Range("M61").Value = -3509933.8
$ws.Range("N61").Value = -2852.5715
$ws.Range("H74").Value = 1146.3
$ws.Range("I74").Value = 1146.3
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1146.3
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -272.3
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 1146.3
$ws.Range("I77").Value = 1146.3
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5731.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1363.5
$ws.Range("N77").ClearContents()
$ws.Range("H117").Value = 31000
$ws.Range("J117").Value = 31000
$ws.Range("L117").Value = 31000
$ws.Range("N117").Value = -40178
$ws.Range("H136").Value = 2565760.5
$ws.Range("I136").Value = 3510145.8
$ws.Range("J136").Value = 2428.5715
$ws.Range("K136").Value = 10530437.4
$ws.Range("L136").Value = 7285.7145
$ws.Range("M136").Value = -10527887.4
$ws.Range("N136").Value = -12385.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 748.5333000000001
$ws.Range("I58").Value = 742.75
$ws.Range("J58").Value = 829.5
$ws.Range("K58").Value = 742.75
$ws.Range("L58").Value = 829.5
$ws.Range("M58").Value = -539.75
$ws.Range("N58").Value = -1235.5
$ws.Range("H132").Value = 28203.105
$ws.Range("I132").Value = 1397.3
$ws.Range("J132").Value = 128724.875
$ws.Range("K132").Value = 4191.9
$ws.Range("L132").Value = 386174.625
$ws.Range("M132").Value = -1661.9
$ws.Range("N132").Value = -391234.625
$ws.Range("H134").Value = 1823.125
$ws.Range("I134").Value = 1325.4849
$ws.Range("J134").Value = 4169.143
$ws.Range("K134").Value = 3976.4547
$ws.Range("L134").Value = 12507.429
$ws.Range("M134").Value = -1441.4547
$ws.Range("N134").Value = -17577.429
$ws.Range("H136").Value = 748.5333000000001
$ws.Range("I136").Value = 742.75
$ws.Range("J136").Value = 829.5
$ws.Range("K136").Value = 2228.25
$ws.Range("L136").Value = 2488.5
$ws.Range("M136").Value = 321.75
$ws.Range("N136").Value = -7588.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 1334.2916
$ws.Range("I38").Value = 1129.875
$ws.Range("J38").Value = 1436.5
$ws.Range("K38").Value = 3389.625
$ws.Range("L38").Value = 4309.5
$ws.Range("M38").Value = -3042.625
$ws.Range("N38").Value = -5003.5
$ws.Range("H80").Value = 2061.818
$ws.Range("I80").Value = 950
$ws.Range("J80").Value = 2308.889
$ws.Range("K80").Value = 2850
$ws.Range("L80").Value = 6926.667
$ws.Range("M80").Value = -1914
$ws.Range("N80").Value = -8798.667000000001
$ws.Range("H83").Value = 2061.818
$ws.Range("I83").Value = 950
$ws.Range("J83").Value = 2308.889
$ws.Range("K83").Value = 8550
$ws.Range("L83").Value = 20780.001
$ws.Range("M83").Value = -3870
$ws.Range("N83").Value = -30140.001
$ws.Range("H131").Value = 20835044
$ws.Range("J131").Value = 33334818
$ws.Range("L131").Value = 100004454
$ws.Range("N131").Value = -100014534

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1545.1111
$ws.Range("I102").Value = 1124.5
$ws.Range("J102").Value = 3395.8
$ws.Range("K102").Value = 1124.5
$ws.Range("L102").Value = 3395.8
$ws.Range("M102").Value = 497.5
$ws.Range("N102").Value = -6639.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 19644.719
$ws.Range("I132").Value = 26067.096
$ws.Range("J132").Value = 1662.0667
$ws.Range("K132").Value = 78201.288
$ws.Range("L132").Value = 4986.2001
$ws.Range("M132").Value = -75671.288
$ws.Range("N132").Value = -10046.2001
$ws.Range("H136").Value = 5319.8057
$ws.Range("I136").Value = 4998.154
$ws.Range("J136").Value = 6156.1
$ws.Range("K136").Value = 14994.462
$ws.Range("L136").Value = 18468.3
$ws.Range("M136").Value = -12444.462
$ws.Range("N136").Value = -23568.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1831.3684
$ws.Range("I122").Value = 1815.4615
$ws.Range("J122").Value = 1865.8334
$ws.Range("K122").Value = 5446.3845
$ws.Range("L122").Value = 5597.5002
$ws.Range("M122").Value = -2996.3845
$ws.Range("N122").Value = -10497.5002
$ws.Range("H126").Value = 33334560
$ws.Range("I126").Value = 66667776
$ws.Range("J126").Value = 1345.3334
$ws.Range("K126").Value = 200003328
$ws.Range("L126").Value = 4036.0002
$ws.Range("M126").Value = -200000858
$ws.Range("N126").Value = -8976.0002
$ws.Range("H132").Value = 57310860
$ws.Range("I132").Value = 90401520
$ws.Range("J132").Value = 2159767
$ws.Range("K132").Value = 271204560
$ws.Range("L132").Value = 6479301
$ws.Range("M132").Value = -271202030
$ws.Range("N132").Value = -6484361
$ws.Range("H136").Value = 24407.117
$ws.Range("I136").Value = 30967.303
$ws.Range("J136").Value = 2758.5
$ws.Range("K136").Value = 92901.909
$ws.Range("L136").Value = 8275.5
$ws.Range("M136").Value = -90351.909
$ws.Range("N136").Value = -13375.5

Write-Host "Applied all edits"